# Repull data, push all data, mean calculation
# Update column F (dSF) values per row to reflect refreshed source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -5
    4  = -5
    5  = 4
    7  = -1
    9  = -5
    10 = -1
    11 = -6
    12 = -2
    13 = -1
    14 = -4
    16 = 7
    17 = -1
    18 = 3
    19 = 4
    20 = 2
    21 = -1
    22 = 2
    23 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
